# Auto-generated Excel COM-interop edit script
# Applies the cryptos.xlsx price/volume update described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.801.94"
$ws.Range("E2").Value = "  -4.10%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.129.45"
$ws.Range("E3").Value = "  -3.71%  "

# Row 4
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.51"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.88"
$ws.Range("E6").Value = "  -7.26%  "

# Row 7
$ws.Range("E7").Value = "  -0.15%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.131.70"
$ws.Range("E8").Value = "  -3.61%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  -4.62%  "

# Row 10
$ws.Range("E10").Value = "  -7.42%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.38"
$ws.Range("E11").Value = "  -8.15%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.470"
$ws.Range("E12").Value = "  -6.14%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  -8.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.27"
$ws.Range("E14").Value = "  -9.80%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.646.00"
$ws.Range("E15").Value = "  -3.89%  "

# Row 16
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.114"
$ws.Range("E16").Value = "  +0.73%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.871.36"
$ws.Range("E17").Value = "  -4.22%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.137.76"
$ws.Range("E18").Value = "  -3.67%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.86"
$ws.Range("E19").Value = "  -8.12%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.97"
$ws.Range("E20").Value = "  -6.03%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.71"
$ws.Range("E21").Value = "  -4.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.701"
$ws.Range("E22").Value = "  -6.71%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.76"
$ws.Range("E23").Value = "  -4.63%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.55"
$ws.Range("E24").Value = "  -8.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.76"
$ws.Range("E25").Value = "  -4.98%  "

# Row 26
$ws.Range("E26").Value = "  -0.10%  "

# Row 27
$ws.Range("E27").Value = "  -5.36%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.38"
$ws.Range("E28").Value = "  -8.52%  "

# Row 29
$ws.Range("E29").Value = "  -9.72%  "

# Row 30
$ws.Range("E30").Value = "  -3.26%  "

# Row 31
$ws.Range("E31").Value = "  -14.58%  "

# Row 32
$ws.Range("E32").Value = "  -6.37%  "

# Row 33
$ws.Range("E33").Value = "  -0.12%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.05"
$ws.Range("E34").Value = "  -7.03%  "

# Row 35
$ws.Range("E35").Value = "  -4.30%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.95"
$ws.Range("E36").Value = "  -8.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.50"
$ws.Range("E37").Value = "  -3.85%  "

# Row 38
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0730"
$ws.Range("E38").Value = "  -6.69%  "

# Row 39
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "459.21"
$ws.Range("E39").Value = "  -7.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("E40").Value = "  -13.52%  "

# Row 41
$ws.Range("E41").Value = "  -7.73%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.39"
$ws.Range("E42").Value = "  -4.99%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.118"
$ws.Range("E43").Value = "  -8.12%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.836.09"
$ws.Range("E44").Value = "  -5.06%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.264"
$ws.Range("E45").Value = "  -10.25%  "

# Row 46
$ws.Range("E46").Value = "  -10.82%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.32"
$ws.Range("E47").Value = "  -8.93%  "

# Row 48
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  +0.03%  "

# Row 49
$ws.Range("E49").Value = "  -8.11%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.114"
$ws.Range("E50").Value = "  -4.74%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.08"
$ws.Range("E51").Value = "  -1.60%  "

